$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: temporarily unhide every row whose RunMode value will change, so the
#     engine does not collapse its row height while writing the new cell value. ---
$ws.Range("A84:A87").EntireRow.Hidden = $false
$ws.Range("A90:A109").EntireRow.Hidden = $false
$ws.Range("A125:A165").EntireRow.Hidden = $false
$ws.Range("A207:A210").EntireRow.Hidden = $false
$ws.Range("A213:A232").EntireRow.Hidden = $false
$ws.Range("A330:A333").EntireRow.Hidden = $false
$ws.Range("A336:A355").EntireRow.Hidden = $false
$ws.Range("A365:A411").EntireRow.Hidden = $false
$ws.Range("A453:A493").EntireRow.Hidden = $false

# --- Step 2: set RunMode (column C) values while all affected rows are visible ---
$ws.Range("C84:C87").Value2 = "Yes"
$ws.Range("C90:C109").Value2 = "Yes"
$ws.Range("C125:C128").Value2 = "Yes"
$ws.Range("C131:C150").Value2 = "Yes"
$ws.Range("C207:C210").Value2 = "Yes"
$ws.Range("C213:C232").Value2 = "Yes"
$ws.Range("C330:C333").Value2 = "Yes"
$ws.Range("C336:C355").Value2 = "Yes"
$ws.Range("C371:C374").Value2 = "Yes"
$ws.Range("C377:C396").Value2 = "Yes"
$ws.Range("C453:C456").Value2 = "Yes"
$ws.Range("C459:C478").Value2 = "Yes"

$ws.Range("C129:C130").Value2 = "No"
$ws.Range("C151:C165").Value2 = "No"
$ws.Range("C365:C370").Value2 = "No"
$ws.Range("C375:C376").Value2 = "No"
$ws.Range("C397:C411").Value2 = "No"
$ws.Range("C457:C458").Value2 = "No"
$ws.Range("C479:C493").Value2 = "No"

# --- Step 3: apply the final Hidden state now that values are written ---
$ws.Range("A85:A87").EntireRow.Hidden = $true
$ws.Range("A95:A109").EntireRow.Hidden = $true
$ws.Range("A126:A130").EntireRow.Hidden = $true
$ws.Range("A136:A165").EntireRow.Hidden = $true
$ws.Range("A208:A210").EntireRow.Hidden = $true
$ws.Range("A218:A232").EntireRow.Hidden = $true
$ws.Range("A331:A333").EntireRow.Hidden = $true
$ws.Range("A341:A355").EntireRow.Hidden = $true
$ws.Range("A365:A370").EntireRow.Hidden = $true
$ws.Range("A372:A376").EntireRow.Hidden = $true
$ws.Range("A382:A411").EntireRow.Hidden = $true
$ws.Range("A454:A458").EntireRow.Hidden = $true
$ws.Range("A464:A493").EntireRow.Hidden = $true
# --- Step 4: sheet view (selection + scroll position) ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 462
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K462").Select()

# --- Step 5: autoFilter changes ---
# Turn off the old "Leave Probation Period according to Custom Months" (column G) filter
# and move/replace it with "Pro rata" (column K) = No; flip "After Probation period"
# (column M) filter from Yes to No.
$rngAF = $ws.Range("A1:Y493")
$rngAF.AutoFilter(7)
$rngAF.AutoFilter(11, @("No"), 7)
$rngAF.AutoFilter(13, @("No"), 7)
